# SFUN_QTR_FIN.xlsx - "Doing Updates for Financials"
# Insert a new quarter column (D) in front of the existing data, shifting
# the previously-existing quarters one column to the right, and populate
# the new column with the latest quarter's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D; this pushes the old D:K data to E:L and
# creates a blank column D (default formatting).
$ws.Columns("D:D").Insert()

# Copy the number formatting/styles from the (now shifted) first data
# column E into the new D column, so the new quarter's cells keep the
# same look (date format on the header row, number format on data rows)
# as the rest of the table.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rows 5, 6, 36, 37, 78 and 79 have no quarterly data at all (they are
# either section headers holding only a label in column A/B, or blank
# spacer rows that don't otherwise exist in the sheet); the paste above
# would have stamped a blank, styled D cell onto them, so remove those
# again.
$ws.Range("D5").Clear()
$ws.Range("D6").Clear()
$ws.Range("D36").Clear()
$ws.Range("D37").Clear()
$ws.Range("D78").Clear()
$ws.Range("D79").Clear()

# ---- Income Statement -------------------------------------------------
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 83600
$ws.Range("D9").Value = 17800
$ws.Range("D10").Value = 65800
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 68700
$ws.Range("D18").Value = 14900
$ws.Range("D20").Value = -3800
$ws.Range("D21").Value = "NA"
$ws.Range("D22").Value = 5800
$ws.Range("D23").Value = 5400
$ws.Range("D24").Value = 3100
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 2300
$ws.Range("D27").Value = 2300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 3800
$ws.Range("D33").Value = 2300
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 2300

# ---- Balance Sheet ------------------------------------------------------
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 151200
$ws.Range("D42").Value = 40200
$ws.Range("D43").Value = 225600
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 245700
$ws.Range("D46").Value = 662700
$ws.Range("D47").Value = 406800
$ws.Range("D48").Value = 737400
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 78500
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1885400

# Row 57 received brand-new source data: the new column gets a real
# value, while the previously-existing quarters (now E:J) came back as
# "NA" from the data feed for this particular line item; K/L keep 0.
$ws.Range("D57").Value = 300
$ws.Range("E57:J57").Value = "NA"

$ws.Range("D58").Value = 280400
$ws.Range("D59").Value = 323200
$ws.Range("D60").Value = 603900
$ws.Range("D61").Value = 417100
$ws.Range("D62").Value = 222300
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1244000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 320500
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 641400
$ws.Range("D77").Value = 0

# ---- Cash Flow ----------------------------------------------------------
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = 2300
$ws.Range("D83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 0
$ws.Range("D91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 0
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 0
